$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "Visualizing Texas: powerpoint_presentation"
#                -> "Texas Housing Prices: powerpoint_presentation"
$title = $s.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# Run 1: "Visualizing" (chars 1-11) -> "Texas"
$runVisualizing = $tr.Characters(1, 11)
$runVisualizing.Text = "Texas"

# Run 3: "Texas:" is now at position 7, length 6 -> becomes "Housing Prices:"
$runTexasColon = $tr.Characters(7, 6)
$runTexasColon.Text = "Housing Prices:"

# Split "Housing Prices:" (pos 7, len 15) into three runs:
#   "Housing" (7,7) + " " (14,1) + "Prices:" (15,7)
$runSpace = $tr.Characters(14, 1)
$runSpace.Text = " "

$runPrices = $tr.Characters(15, 7)
$runPrices.Text = "Prices:"

# --- Remove the "Date Placeholder 3" shape entirely.
# The first Delete() clears the placeholder's content but the layout's
# empty placeholder gets re-materialized in its place (with a new name);
# the second Delete() removes that leftover empty placeholder too.
$s.Shapes.Item(3).Delete()
$s.Shapes.Item(3).Delete()
